$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rules")

# Rules sheet, row 11 ("R40" rule row): the "To" cell (B11) changes from the
# text "R40" to the text "1". A leading apostrophe forces Excel to store the
# numeric-looking entry as text (adds a new shared string "1") instead of
# coercing it to the number 1.
$ws.Range("B11").Value = "'1"
